$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 1519.7273
$ws.Range("I40").Value = 760
$ws.Range("J40").Value = 1953.8572
$ws.Range("K40").Value = 760
$ws.Range("L40").Value = 1953.8572
$ws.Range("M40").Value = -585
$ws.Range("N40").Value = -2303.8572

$ws.Range("H64").Value = 5300
$ws.Range("J64").Value = 5300
$ws.Range("L64").Value = 5300
$ws.Range("N64").Value = -5796

$ws.Range("H67").Value = 5300
$ws.Range("J67").Value = 5300
$ws.Range("L67").Value = 5300
$ws.Range("N67").Value = -7016

$ws.Range("H74").Value = 3468.4375
$ws.Range("I74").Value = 2585
$ws.Range("K74").Value = 2585
$ws.Range("M74").Value = -1649

$ws.Range("H77").Value = 3468.4375
$ws.Range("I77").Value = 2585
$ws.Range("K77").Value = 12925
$ws.Range("M77").Value = -8245

$ws.Range("H80").Value = 10153012
$ws.Range("I80").Value = 611.1667
$ws.Range("J80").Value = 20305412
$ws.Range("K80").Value = 1833.5001
$ws.Range("L80").Value = 60916236
$ws.Range("M80").Value = -835.5001
$ws.Range("N80").Value = -60918232

$ws.Range("H83").Value = 10153012
$ws.Range("I83").Value = 611.1667
$ws.Range("J83").Value = 20305412
$ws.Range("K83").Value = 5500.5003
$ws.Range("L83").Value = 182748708
$ws.Range("M83").Value = -508.5002999999997
$ws.Range("N83").Value = -182758692

$ws.Range("H113").Value = 50004812
$ws.Range("J113").Value = 6142.857
$ws.Range("L113").Value = 6142.857
$ws.Range("N113").Value = -12650.857

$ws.Range("H135").Value = 11908176
$ws.Range("J135").Value = 41676850
$ws.Range("L135").Value = 375091650
$ws.Range("N135").Value = -375096720

$ws.Range("H141").Value = 1229.1177
$ws.Range("I141").Value = 822.2955
$ws.Range("K141").Value = 2466.8865
$ws.Range("M141").Value = 2713.1135

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H102").Value = 1100
$ws.Range("I102").Value = 1100
$ws.Range("K102").Value = 1100
$ws.Range("M102").Value = 522

$ws.Range("H132").Value = 13467.643
$ws.Range("I132").Value = 1571.7354
$ws.Range("J132").Value = 64025.25
$ws.Range("K132").Value = 4715.206200000001
$ws.Range("L132").Value = 192075.75
$ws.Range("M132").Value = -2185.206200000001
$ws.Range("N132").Value = -197135.75

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 3353.4614
$ws.Range("I105").Value = 3611.875
$ws.Range("J105").Value = 2940
$ws.Range("K105").Value = 3611.875
$ws.Range("L105").Value = 2940
$ws.Range("M105").Value = -1864.875
$ws.Range("N105").Value = -6434

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H2").Value = 701.3333
$ws.Range("I2").Value = 104
$ws.Range("J2").Value = 1000
$ws.Range("K2").Value = 104
$ws.Range("L2").Value = 1000
$ws.Range("M2").Value = 9
$ws.Range("N2").Value = -1226

$ws.Range("H3").Value = 14000
$ws.Range("J3").Value = 14000
$ws.Range("L3").Value = 14000
$ws.Range("N3").Value = -14226

$ws.Range("H7").Value = 84.2
$ws.Range("I7").Value = 99.666664
$ws.Range("J7").Value = 61
$ws.Range("K7").Value = 99.666664
$ws.Range("L7").Value = 61
$ws.Range("M7").Value = 13.333336
$ws.Range("N7").Value = -287

$ws.Range("H22").Value = 431.375
$ws.Range("I22").Value = 289.8
$ws.Range("J22").Value = 667.3333
$ws.Range("K22").Value = 289.8
$ws.Range("L22").Value = 667.3333
$ws.Range("M22").Value = 60.19999999999999
$ws.Range("N22").Value = -1367.3333

$ws.Range("H31").Value = 2952.8667
$ws.Range("I31").Value = 1698.0303
$ws.Range("J31").Value = 6403.6665
$ws.Range("K31").Value = 1698.0303
$ws.Range("L31").Value = 6403.6665
$ws.Range("M31").Value = -1403.0303
$ws.Range("N31").Value = -6993.6665

$ws.Range("H34").Value = 2952.8667
$ws.Range("I34").Value = 1698.0303
$ws.Range("J34").Value = 6403.6665
$ws.Range("K34").Value = 1698.0303
$ws.Range("L34").Value = 6403.6665
$ws.Range("M34").Value = -1496.0303
$ws.Range("N34").Value = -6807.6665

$ws.Range("H99").Value = 21742880
$ws.Range("I99").Value = 3580
$ws.Range("J99").Value = 38465416
$ws.Range("K99").Value = 3580
$ws.Range("L99").Value = 38465416
$ws.Range("M99").Value = -2082
$ws.Range("N99").Value = -38468412

$ws.Range("H126").Value = 21742880
$ws.Range("I126").Value = 3580
$ws.Range("J126").Value = 38465416
$ws.Range("K126").Value = 10740
$ws.Range("L126").Value = 115396248
$ws.Range("M126").Value = -8270
$ws.Range("N126").Value = -115401188

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 771.37
$ws.Range("J131").Value = 783.65265
$ws.Range("L131").Value = 2350.95795
$ws.Range("N131").Value = -12430.95795

$ws.Range("H132").Value = 433
$ws.Range("I132").Value = 433
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 3897
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -1367
$ws.Range("N132").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H7").Value = 4600000
$ws.Range("J7").Value = 4000000
$ws.Range("L7").Value = 4000000
$ws.Range("N7").Value = -4000224

$ws.Range("H8").Value = 4600000
$ws.Range("J8").Value = 4000000
$ws.Range("L8").Value = 4000000
$ws.Range("N8").Value = -4000278

$ws.Range("H113").Value = 5250
$ws.Range("I113").Value = 4000
$ws.Range("J113").Value = 5666.6665
$ws.Range("K113").Value = 4000
$ws.Range("L113").Value = 5666.6665
$ws.Range("M113").Value = -1830
$ws.Range("N113").Value = -10006.6665

$ws.Range("H132").Value = 25934.182
$ws.Range("I132").Value = 3590.4736
$ws.Range("J132").Value = 167444.33
$ws.Range("K132").Value = 10771.4208
$ws.Range("L132").Value = 502332.99
$ws.Range("M132").Value = -8241.4208
$ws.Range("N132").Value = -507392.99

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 1326.7142
$ws.Range("I46").Value = 997.4
$ws.Range("J46").Value = 2150
$ws.Range("K46").Value = 997.4
$ws.Range("L46").Value = 2150
$ws.Range("M46").Value = -809.4
$ws.Range("N46").Value = -2526

$ws.Range("H132").Value = 1054.8235
$ws.Range("I132").Value = 1054.8235
$ws.Range("K132").Value = 3164.4705
$ws.Range("M132").Value = -634.4704999999999
